# Tidsregistrering - Dan Toft.xlsx
# "Gammel Tidsregistereing der aldrig blev pushet"
#
# Fills in activity descriptions (column D) and end-times (column F) for
# rows 26-33 on the "Tidsregistrering" sheet, corrects the date in B28
# (typed directly, no longer the shared date formula), and moves the
# active-cell selection to F34.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsregistrering")

# Row 26 - Guikode, ends 15:30
$ws.Range("D26").Value = "Guikode"
$ws.Range("F26").Value = 0.64583333333333337

# Row 27 - Det kan jeg ikke huske..., ends 15:30
$ws.Range("D27").Value = "Det kan jeg ikke huske…"
$ws.Range("F27").Value = 0.64583333333333337

# Row 28 - new day (2017-03-20) typed directly into B28, design oc5, ends 9:30
$ws.Range("B28").Value = 42814
$ws.Range("D28").Value = "design oc5"
$ws.Range("F28").Value = 0.39583333333333331

# Row 29 - Review OC11, ends 9:40
$ws.Range("D29").Value = "Review OC11"
$ws.Range("F29").Value = 0.40277777777777773

# Row 30 - Testsuite oc9, ends 12:25
$ws.Range("D30").Value = "Testsuite oc9"
$ws.Range("F30").Value = 0.51736111111111105

# Row 31 - design review, ends 14:00
$ws.Range("D31").Value = "design review"
$ws.Range("F31").Value = 0.58333333333333337

# Row 32 - Implamenter et eller andet, ends 15:30
$ws.Range("D32").Value = "Implamenter et eller andet"
$ws.Range("F32").Value = 0.64583333333333337

# Row 33 - Review af design på OC13, ends 9:00
$ws.Range("D33").Value = "Review af design på OC13"
$ws.Range("F33").Value = 0.375

# Move the active selection as recorded in the saved sheet view
$ws.Range("F34").Select()
